$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '31.373.39'
$ws.Range('E2').Value = '  +3.52%  '
$ws.Range('D3').Value = '2.006.88'
$ws.Range('E3').Value = '  +7.54%  '
$ws.Range('D4').Value = '0.9994'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '0.7899'
$ws.Range('E5').Value = '  +67.22%  '
$ws.Range('D6').Value = '259.33'
$ws.Range('E6').Value = '  +6.52%  '
$ws.Range('D7').Value = '0.9992'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '0.3582'
$ws.Range('E8').Value = '  +24.82%  '
$ws.Range('D9').Value = '28.44'
$ws.Range('E9').Value = '  +31.50%  '
$ws.Range('E10').Value = '  +8.98%  '
$ws.Range('D11').Value = '0.8492'
$ws.Range('E11').Value = '  +17.91%  '
$ws.Range('D12').Value = '0.08118'
$ws.Range('E12').Value = '  +4.19%  '
$ws.Range('D13').Value = '2.006.55'
$ws.Range('E13').Value = '  +7.55%  '
$ws.Range('D14').Value = '101.01'
$ws.Range('E14').Value = '  +4.10%  '
$ws.Range('D15').Value = '5.606'
$ws.Range('E15').Value = '  +9.04%  '
$ws.Range('D16').Value = '277.28'
$ws.Range('E16').Value = '  -1.09%  '
$ws.Range('B17').Value = 'Avalanche'
$ws.Range('C17').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D17').Value = '14.85'
$ws.Range('E17').Value = '  +14.24%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '31.385.62'
$ws.Range('E18').Value = '  +3.58%  '
$ws.Range('D19').Value = '5.930'
$ws.Range('E19').Value = '  +13.45%  '
$ws.Range('D20').Value = '0.000007925'
$ws.Range('E20').Value = '  +6.26%  '
$ws.Range('D21').Value = '2.270.73'
$ws.Range('E21').Value = '  +7.84%  '
$ws.Range('D22').Value = '0.9993'
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').Value = '0.9998'
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').Value = '7.191'
$ws.Range('E24').Value = '  +15.06%  '
$ws.Range('E25').Value = '  +12.20%  '
$ws.Range('D26').Value = '0.1508'
$ws.Range('E26').Value = '  +56.50%  '
$ws.Range('D27').Value = '164.01'
$ws.Range('E27').Value = '  +1.18%  '
$ws.Range('D28').Value = '20.01'
$ws.Range('E28').Value = '  +7.10%  '
$ws.Range('D29').Value = '2.376'
$ws.Range('E29').Value = '  +26.64%  '
$ws.Range('E30').Value = '  +10.06%  '
$ws.Range('D31').Value = '4.630'
$ws.Range('E31').Value = '  +9.82%  '
$ws.Range('D32').Value = '1.357'
$ws.Range('E32').Value = '  +3.30%  '
$ws.Range('D33').Value = '4.406'
$ws.Range('E33').Value = '  +7.20%  '
$ws.Range('D34').Value = '0.05238'
$ws.Range('E34').Value = '  +9.71%  '
$ws.Range('D36').Value = '0.7652'
$ws.Range('E36').Value = '  +11.99%  '
$ws.Range('D37').Value = '2.806'
$ws.Range('E37').Value = '  +3.55%  '
$ws.Range('D38').Value = '0.02015'
$ws.Range('D39').Value = '2.948'
$ws.Range('E39').Value = '  +3.98%  '
$ws.Range('D40').Value = '80.53'
$ws.Range('E40').Value = '  +7.11%  '
$ws.Range('D41').Value = '6.689'
$ws.Range('E41').Value = '  +7.94%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = '0.4752'
$ws.Range('E42').Value = '  +13.03%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '2.174'
$ws.Range('E43').Value = '  +12.29%  '
$ws.Range('D44').Value = '0.8600'
$ws.Range('E44').Value = '  +4.05%  '
$ws.Range('D45').Value = '105.40'
$ws.Range('E45').Value = '  +4.91%  '
$ws.Range('D46').Value = '0.9995'
$ws.Range('E46').Value = '  +0.07%  '
$ws.Range('D47').Value = '7.743'
$ws.Range('E47').Value = '  +11.40%  '
$ws.Range('D48').Value = '9.976'
$ws.Range('E48').Value = '  +3.97%  '
$ws.Range('D49').Value = '0.4370'
$ws.Range('E49').Value = '  +12.76%  '
$ws.Range('D50').Value = '36.89'
$ws.Range('E50').Value = '  +5.59%  '
$ws.Range('D51').Value = '0.1192'
$ws.Range('E51').Value = '  +14.85%  '
